$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("Data")
$codebook = $wb.Worksheets.Item("Codebook")

# ---------------------------------------------------------------------------
# 1. Data sheet: add "Income" (D) and "BMI" (E) columns
# ---------------------------------------------------------------------------

# Headers (bold, matching existing header style)
$data.Range("D1").Value = "Income"
$data.Range("E1").Value = "BMI"
$data.Range("D1:E1").Font.Bold = $true

# Income values (column D)
$incomeValues = @(70000, 80000, 90000, 95000, 75000, 85000, 89000, 85000, 92000, 78000, 110000, 100000, 60000, 96000)
for ($i = 0; $i -lt $incomeValues.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 4).Value = $incomeValues[$i]
}

# Apply currency number format to the Income data cells
$data.Range("D2:D15").NumberFormat = "$#,##0_);[Red]($#,##0)"

# BMI values (column E) - static entered values for rows 2-14
$bmiValues = @(24, 22, 22, 23, 24, 22, 36, 39, 22, 22, 22, 25, 19)
for ($i = 0; $i -lt $bmiValues.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 5).Value = $bmiValues[$i]
}

# Row 15 BMI is computed with a formula
$data.Range("E15").Formula = "=ROUND(B15/((A15/100)^2),0)"

# Restore/adjust selection on the Data sheet (it is no longer the active tab)
[void]$data.Activate()
[void]$data.Range("M14").Select()

# ---------------------------------------------------------------------------
# 2. Codebook sheet: document the two new variables
# ---------------------------------------------------------------------------

$codebook.Range("A5").Value = "Income"
$codebook.Range("B5").Value = "Individual Income in usd"
$codebook.Range("C5").Value = "numeric value >0 or NA"

$codebook.Range("A6").Value = "BMI"
$codebook.Range("B6").Value = "Body mass index in kg/m2"
$codebook.Range("C6").Value = "numeric value >0 or NA"

# Widen column B to fit the new, longer descriptions
$codebook.Columns.Item(2).AutoFit()

# Codebook becomes the active sheet/selection
[void]$codebook.Activate()
[void]$codebook.Range("B7").Select()
